$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update actor name from "Diretor" to "Diretor do Hotel"
$ws.Range("A2").Value = "Diretor do Hotel"

# Update the selected range to match the new view state (A7:A10, active cell A7)
$ws.Activate()
$ws.Range("A7:A10").Select()
